$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (rows 2-51) to remain Text so numeric-looking
# strings (e.g. "413.10", "0.0000224") are not coerced into numbers,
# matching the original inline-string cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

# Apply the updated crypto price/volume snapshot (and the two name/link
# swaps: Uniswap<->WrappedEther at rows 17-18, EnergySwap<->PEPE at rows 49-50)
$ws.Range('D2').Value = '61.830.94'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '3.412.23'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '413.10'
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('D6').Value = '129.03'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  -2.75%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').Value = '0.141'
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('D11').Value = '42.60'
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('D12').Value = '0.0000224'
$ws.Range('E12').Value = '  +5.38%  '
$ws.Range('E13').Value = '  +2.05%  '
$ws.Range('D14').Value = '3.951.53'
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('E16').Value = '  -1.46%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '12.74'
$ws.Range('E17').Value = '  +5.23%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.397.49'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').Value = '61.880.50'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').Value = '478.88'
$ws.Range('E21').Value = '  +8.37%  '
$ws.Range('D22').Value = '91.10'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').Value = '3.26'
$ws.Range('E23').Value = '  +3.40%  '
$ws.Range('D24').Value = '13.06'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').Value = '3.30'
$ws.Range('E26').Value = '  +11.21%  '
$ws.Range('D27').Value = '33.20'
$ws.Range('E27').Value = '  -1.86%  '
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('D29').Value = '7.67'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('D30').Value = '11.81'
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('E31').Value = '  -1.35%  '
$ws.Range('D32').Value = '0.166'
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('E33').Value = '  -2.62%  '
$ws.Range('D34').Value = '40.86'
$ws.Range('E34').Value = '  -3.69%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').Value = '58.47'
$ws.Range('E36').Value = '  +8.85%  '
$ws.Range('D37').Value = '0.0485'
$ws.Range('E37').Value = '  -3.36%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '3.01'
$ws.Range('E39').Value = '  +4.01%  '
$ws.Range('D40').Value = '0.325'
$ws.Range('E40').Value = '  +3.73%  '
$ws.Range('D41').Value = '148.50'
$ws.Range('E41').Value = '  +5.22%  '
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('E44').Value = '  +5.97%  '
$ws.Range('E45').Value = '  +7.22%  '
$ws.Range('D46').Value = '4.21'
$ws.Range('E46').Value = '  +3.13%  '
$ws.Range('D47').Value = '2.33'
$ws.Range('E47').Value = '  +17.80%  '
$ws.Range('D48').Value = '16.34'
$ws.Range('E48').Value = '  -1.16%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '22.27'
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.0₃0529'
$ws.Range('E50').Value = '  +24.48%  '
$ws.Range('D51').Value = '113.22'
$ws.Range('E51').Value = '  +7.46%  '
